# Rename the header row strings to reflect the new "FV2310/FV2404" naming
# scheme instead of "old/new" (columns A:J -> _FV2310, columns L:U -> _FV2404;
# column K is just "diff" and is left untouched).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value = ($cell.Value2 -replace "_old$", "_FV2310")
}
for ($c = 12; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value = ($cell.Value2 -replace "_new$", "_FV2404")
}

# Turn the data range into an Excel Table ("Table1") so it gets filter
# buttons / structured references, matching the header names above.
$rng = $ws.Range("A1:U62")
$lo = $ws.ListObjects.Add(1, $rng, $null, 1)
$lo.Name = "Table1"

# Freeze the header row (split below row 1, then lock it).
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
# Leave the selection back on A1 (matches the default/untouched selection
# state used by the exporter instead of lingering on the freeze anchor).
$ws.Range("A1").Select() | Out-Null
